## "Add all Consumer Widget !!"
# Rework the CONSUMER sheet so it carries a new "Wifi" component (SKY85207 / PA)
# alongside the existing SoC info, add the voltage/current input rows, and
# switch the active window focus over to the CONSUMER tab.

$wb = $excel.ActiveWorkbook

$wsDCDC     = $wb.Worksheets.Item("DCDC")
$wsPSU      = $wb.Worksheets.Item("PSU")
$wsCONSUMER = $wb.Worksheets.Item("CONSUMER")

# --- CONSUMER sheet: new headers + new Wifi column + new rows -------------

$wsCONSUMER.Range("A1").Value = "Name"
$wsCONSUMER.Range("B1").Value = "Soc"
$wsCONSUMER.Range("C1").Value = "Wifi"

$wsCONSUMER.Range("A2").Value = "Ref"
$wsCONSUMER.Range("B2").Value = "BCM72180"
$wsCONSUMER.Range("C2").Value = "SKY85207"

$wsCONSUMER.Range("A3").Value = "Info"
$wsCONSUMER.Range("B3").Value = "Avs Core"
$wsCONSUMER.Range("C3").Value = "PA"

$wsCONSUMER.Range("A4").Value = "equivalence code"
$wsCONSUMER.Range("B4").Value = 191477712
$wsCONSUMER.Range("C4").Value = 191479021

$wsCONSUMER.Range("A5").Value = "voltage input"
$wsCONSUMER.Range("B5").Value = 0.8
$wsCONSUMER.Range("C5").Value = 3.3

$wsCONSUMER.Range("A6").Value = "current input"
$wsCONSUMER.Range("B6").Value = 1500
$wsCONSUMER.Range("C6").Value = 250

# Widen column A to fit the new "equivalence code" / "voltage input" labels
# (closest the host's character-width quantization can land to 16.43 chars).
$wsCONSUMER.Columns.Item(1).ColumnWidth = 15.6666666666667

# --- Window / selection state ----------------------------------------------

# DCDC keeps its data untouched, just the remembered selection moves.
$wsDCDC.Range("P26").Select()

# PSU is no longer the focused tab; park its remembered selection at A3.
$wsPSU.Range("A3").Select()

# CONSUMER becomes the active tab, with F11 as the remembered selection.
$wsCONSUMER.Activate()
$wsCONSUMER.Range("F11").Select()
